# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D2").Value = '66.531.64'
$ws.Range("E2").Value = '  +2.64%  '
$ws.Range("D3").Value = '3.201.41'
$ws.Range("E3").Value = '  +1.78%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.200.06'
$ws.Range("E8").Value = '  +1.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("E12").Value = '  +4.11%  '
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.07'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.71%  '
$ws.Range("D15").Value = '3.726.56'
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").Value = '66.508.92'
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("E17").Value = '  +5.49%  '
$ws.Range("D18").Value = '3.196.60'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '514.24'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.41'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.741'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.10'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.10'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.81'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.52%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.31'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.48%  '
$ws.Range("E29").Value = '  +6.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.18'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +16.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.95'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +5.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.32'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("E33").Value = '  +3.87%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '515.13'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +8.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.88'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0900'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.89'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.32%  '
$ws.Range("E41").Value = '  +6.62%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.68%  '
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.303'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +7.76%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0676'
$ws.Range("E44").Value = '  +17.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.98%  '
$ws.Range("D46").Value = '2.928.16'
$ws.Range("E46").Value = '  -2.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.74'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.61%  '
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("E50").Value = '  +4.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.62'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +10.39%  '
